$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "2024-02-22 11:35:05"
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0.001
$ws.Range("J25").Value = 0.05
$ws.Range("K25").Value = 0.003
$ws.Range("L25").Value = 100
$ws.Range("M25").Value = 500
$ws.Range("N25").Value = 10
$ws.Range("O25").Value = 5
$ws.Range("P25").Value = 3
$ws.Range("Q25").Value = 500
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 100
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = "Data/bombayauto.xlsx"
$ws.Range("W25").Value = 98000
$ws.Range("A26").Value = "2024-02-24 00:03:34"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0.001
$ws.Range("J26").Value = 0.05
$ws.Range("K26").Value = 0.003
$ws.Range("L26").Value = 100
$ws.Range("M26").Value = 500
$ws.Range("N26").Value = 10
$ws.Range("O26").Value = 5
$ws.Range("P26").Value = 2
$ws.Range("Q26").Value = 200
$ws.Range("R26").Value = 8
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 100
$ws.Range("U26").Value = 1
$ws.Range("V26").Value = "Data/bombayauto.xlsx"
$ws.Range("W26").Value = -30800
$ws.Range("A27").Value = "2024-02-24 00:20:11"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0.001
$ws.Range("J27").Value = 0.05
$ws.Range("K27").Value = 0.003
$ws.Range("L27").Value = 100
$ws.Range("M27").Value = 500
$ws.Range("N27").Value = 10
$ws.Range("O27").Value = 5
$ws.Range("P27").Value = 2
$ws.Range("Q27").Value = 500
$ws.Range("R27").Value = 8
$ws.Range("S27").Value = 3
$ws.Range("T27").Value = 100
$ws.Range("U27").Value = 1
$ws.Range("V27").Value = "Data/bombayauto.xlsx"
$ws.Range("W27").Value = 53000
$ws.Range("A28").Value = "2024-02-24 01:22:17"
$ws.Range("B28").Value = 43
$ws.Range("C28").Value = 29
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 12
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0.001
$ws.Range("J28").Value = 0.05
$ws.Range("K28").Value = 0.003
$ws.Range("L28").Value = 100
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = 10
$ws.Range("O28").Value = 5
$ws.Range("P28").Value = 2
$ws.Range("Q28").Value = 200
$ws.Range("R28").Value = 8
$ws.Range("S28").Value = 3
$ws.Range("T28").Value = 80
$ws.Range("U28").Value = 0.6744186046511628
$ws.Range("V28").Value = "Data/bombay1.xlsx"
$ws.Range("W28").Value = -160200
$ws.Range("A29").Value = "2024-02-24 13:52:25"
$ws.Range("B29").Value = 14
$ws.Range("C29").Value = 10
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0.001
$ws.Range("J29").Value = 0.05
$ws.Range("K29").Value = 0.003
$ws.Range("L29").Value = 100
$ws.Range("M29").Value = 500
$ws.Range("N29").Value = 10
$ws.Range("O29").Value = 5
$ws.Range("P29").Value = 2
$ws.Range("Q29").Value = 200
$ws.Range("R29").Value = 8
$ws.Range("S29").Value = 3
$ws.Range("T29").Value = 70
$ws.Range("U29").Value = 0.7142857142857143
$ws.Range("V29").Value = "Data/bombayauto.xlsx"
$ws.Range("W29").Value = 43800
